$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: A2 is a plain formula, A3:A20 keep the existing shared-formula
# group (now stepping by +2 instead of +100), and A21 extends the series.
$ws.Range("A2").Formula = "=+A1+2"
$ws.Range("A3:A20").Formula = "=+A2+2"
$ws.Range("A21").Formula = "=+A20+2"

# Column B: replace all formulas/values with plain literal values (triangle wave 0..10..0)
# and drop the "Zarez" (comma) number-format style that the cells used to carry.
$bVals = @(0, 2, 4, 6, 8, 10, 8, 6, 4, 2, 0, 2, 4, 6, 8, 10, 8, 6, 4, 2, 0)
for ($i = 0; $i -lt $bVals.Length; $i++) {
    $cell = $ws.Cells.Item($i + 1, 2)
    $cell.Value = $bVals[$i]
    $cell.Style = "Normalno"
}

# Selection moves from the old C-column leftover to the new A3:A21 range
$ws.Range("A3:A21").Select() | Out-Null
